# "Generate Report for Handback"
# Updates the localization-status workbook to reflect a completed handback:
#  - Status text flips from "Ready for handoff" to "Handed back: in sync with en-US"
#  - zh-cn / de-de detail sheets get their "Latest Target File" / "Latest Handback
#    File" / "Latest Handback DateTime" columns populated for both rows
#  - New hyperlinks are added on the "Latest Target File" cells
#  - A couple of columns widen to fit the newly-written text

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This text is shared by the Overview sheet (zh-cn/de-de status columns)
#    and by the per-language sheets' Status column, so a blanket replace
#    across the workbook keeps every occurrence in sync.
# ---------------------------------------------------------------------------
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")
}

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: fill in the generated handback artifacts for both rows.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

$zhcn.Range("I2").Value = "a.md"
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-30 10:42:24"

$zhcn.Range("I3").Value = "a.md"
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-30 10:42:24"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d65837becd5392bee9364fad3e26e19c9656c20f/e2e/a.md", "", "", "a.md")
$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = 15570276

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d65837becd5392bee9364fad3e26e19c9656c20f/e2e/a.md", "", "", "a.md")
$zhcn.Range("I3").Font.Underline = 2
$zhcn.Range("I3").Font.Color = 15570276

# ---------------------------------------------------------------------------
# 3. de-de sheet: same shape of update, but with its own xlf/timestamp.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664

$dede.Range("I2").Value = "a.md"
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-08-30 10:42:30"

$dede.Range("I3").Value = "a.md"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = "2016-08-30 10:42:30"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d65837becd5392bee9364fad3e26e19c9656c20f/e2e/a.md", "", "", "a.md")
$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = 15570276

$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d65837becd5392bee9364fad3e26e19c9656c20f/e2e/a.md", "", "", "a.md")
$dede.Range("I3").Font.Underline = 2
$dede.Range("I3").Font.Color = 15570276
